$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: re-order match rows (columns F:V) within same-day blocks ---
# (kickoff-time re-sort shuffled which physical row holds which fixture;
#  column A..E - index/country/league/season/date - are untouched)
$data = @{}
$data[3] = $ws.Range("F3:V3").Value2
$data[4] = $ws.Range("F4:V4").Value2
$data[5] = $ws.Range("F5:V5").Value2
$data[6] = $ws.Range("F6:V6").Value2
$data[7] = $ws.Range("F7:V7").Value2
$data[13] = $ws.Range("F13:V13").Value2
$data[14] = $ws.Range("F14:V14").Value2
$data[15] = $ws.Range("F15:V15").Value2
$data[16] = $ws.Range("F16:V16").Value2
$data[21] = $ws.Range("F21:V21").Value2
$data[22] = $ws.Range("F22:V22").Value2
$data[23] = $ws.Range("F23:V23").Value2
$data[24] = $ws.Range("F24:V24").Value2
$data[25] = $ws.Range("F25:V25").Value2
$data[30] = $ws.Range("F30:V30").Value2
$data[31] = $ws.Range("F31:V31").Value2
$data[33] = $ws.Range("F33:V33").Value2
$data[34] = $ws.Range("F34:V34").Value2
$data[39] = $ws.Range("F39:V39").Value2
$data[40] = $ws.Range("F40:V40").Value2
$data[42] = $ws.Range("F42:V42").Value2
$data[43] = $ws.Range("F43:V43").Value2
$data[48] = $ws.Range("F48:V48").Value2
$data[50] = $ws.Range("F50:V50").Value2
$data[57] = $ws.Range("F57:V57").Value2
$data[58] = $ws.Range("F58:V58").Value2
$data[59] = $ws.Range("F59:V59").Value2
$data[60] = $ws.Range("F60:V60").Value2

$ws.Range("F3:V3").Value2 = $data[4]
$ws.Range("F4:V4").Value2 = $data[7]
$ws.Range("F5:V5").Value2 = $data[6]
$ws.Range("F6:V6").Value2 = $data[5]
$ws.Range("F7:V7").Value2 = $data[3]
$ws.Range("F13:V13").Value2 = $data[16]
$ws.Range("F14:V14").Value2 = $data[15]
$ws.Range("F15:V15").Value2 = $data[14]
$ws.Range("F16:V16").Value2 = $data[13]
$ws.Range("F21:V21").Value2 = $data[25]
$ws.Range("F22:V22").Value2 = $data[24]
$ws.Range("F23:V23").Value2 = $data[21]
$ws.Range("F24:V24").Value2 = $data[22]
$ws.Range("F25:V25").Value2 = $data[23]
$ws.Range("F30:V30").Value2 = $data[34]
$ws.Range("F31:V31").Value2 = $data[33]
$ws.Range("F33:V33").Value2 = $data[30]
$ws.Range("F34:V34").Value2 = $data[31]
$ws.Range("F39:V39").Value2 = $data[42]
$ws.Range("F40:V40").Value2 = $data[43]
$ws.Range("F42:V42").Value2 = $data[40]
$ws.Range("F43:V43").Value2 = $data[39]
$ws.Range("F48:V48").Value2 = $data[50]
$ws.Range("F50:V50").Value2 = $data[48]
$ws.Range("F57:V57").Value2 = $data[60]
$ws.Range("F58:V58").Value2 = $data[57]
$ws.Range("F59:V59").Value2 = $data[58]
$ws.Range("F60:V60").Value2 = $data[59]

# --- Step 2: append new rows 74:82 (matches played 10/10 - 29/10) ---
$ws.Range("A73:V73").Copy()
$ws.Range("A74:V82").PasteSpecial(-4122)

$newRows = New-Object 'object[,]' 9,22
$newRows[0,0] = 73
$newRows[0,1] = 'germany'
$newRows[0,2] = 'bundesliga'
$newRows[0,3] = '2023-2024'
$newRows[0,4] = 45226.85416666666
$newRows[0,5] = 'Bochum'
$newRows[0,6] = 2
$newRows[0,7] = 'Mainz'
$newRows[0,8] = 2
$newRows[0,9] = 2.38
$newRows[0,10] = '10/10/2023 14:02'
$newRows[0,11] = 2.85
$newRows[0,12] = '27/10/2023 20:29'
$newRows[0,13] = 3.48
$newRows[0,14] = '10/10/2023 14:02'
$newRows[0,15] = 3.52
$newRows[0,16] = '27/10/2023 20:29'
$newRows[0,17] = 3.06
$newRows[0,18] = '10/10/2023 14:02'
$newRows[0,19] = 2.55
$newRows[0,20] = '27/10/2023 20:29'
$newRows[0,21] = 'https://www.betexplorer.com/football/germany/bundesliga/bochum-mainz/IDt4N0W2/'

$newRows[1,0] = 74
$newRows[1,1] = 'germany'
$newRows[1,2] = 'bundesliga'
$newRows[1,3] = '2023-2024'
$newRows[1,4] = 45227.64583333334
$newRows[1,5] = 'Augsburg'
$newRows[1,6] = 3
$newRows[1,7] = 'Wolfsburg'
$newRows[1,8] = 2
$newRows[1,9] = 2.72
$newRows[1,10] = '10/10/2023 14:02'
$newRows[1,11] = 2.54
$newRows[1,12] = '28/10/2023 14:56'
$newRows[1,13] = 3.62
$newRows[1,14] = '10/10/2023 14:02'
$newRows[1,15] = 3.62
$newRows[1,16] = '28/10/2023 14:56'
$newRows[1,17] = 2.44
$newRows[1,18] = '10/10/2023 14:02'
$newRows[1,19] = 2.8
$newRows[1,20] = '28/10/2023 14:56'
$newRows[1,21] = 'https://www.betexplorer.com/football/germany/bundesliga/augsburg-wolfsburg/d4u8MKo9/'

$newRows[2,0] = 75
$newRows[2,1] = 'germany'
$newRows[2,2] = 'bundesliga'
$newRows[2,3] = '2023-2024'
$newRows[2,4] = 45227.64583333334
$newRows[2,5] = 'Bayern Munich'
$newRows[2,6] = 8
$newRows[2,7] = 'Darmstadt'
$newRows[2,8] = 0
$newRows[2,9] = 1.07
$newRows[2,10] = '10/10/2023 14:25'
$newRows[2,11] = 1.06
$newRows[2,12] = '28/10/2023 15:16'
$newRows[2,13] = 12.93
$newRows[2,14] = '10/10/2023 14:25'
$newRows[2,15] = 17.61
$newRows[2,16] = '28/10/2023 15:29'
$newRows[2,17] = 19.15
$newRows[2,18] = '10/10/2023 14:25'
$newRows[2,19] = 32.15
$newRows[2,20] = '28/10/2023 15:29'
$newRows[2,21] = 'https://www.betexplorer.com/football/germany/bundesliga/bayern-munich-darmstadt/p2wxKRPA/'

$newRows[3,0] = 76
$newRows[3,1] = 'germany'
$newRows[3,2] = 'bundesliga'
$newRows[3,3] = '2023-2024'
$newRows[3,4] = 45227.64583333334
$newRows[3,5] = 'B. Monchengladbach'
$newRows[3,6] = 2
$newRows[3,7] = 'Heidenheim'
$newRows[3,8] = 1
$newRows[3,9] = 1.93
$newRows[3,10] = '10/10/2023 14:28'
$newRows[3,11] = 1.69
$newRows[3,12] = '28/10/2023 15:28'
$newRows[3,13] = 3.9
$newRows[3,14] = '10/10/2023 14:28'
$newRows[3,15] = 4.44
$newRows[3,16] = '28/10/2023 15:29'
$newRows[3,17] = 3.55
$newRows[3,18] = '10/10/2023 14:28'
$newRows[3,19] = 4.7
$newRows[3,20] = '28/10/2023 15:29'
$newRows[3,21] = 'https://www.betexplorer.com/football/germany/bundesliga/b-monchengladbach-heidenheim/j7lIacvd/'

$newRows[4,0] = 77
$newRows[4,1] = 'germany'
$newRows[4,2] = 'bundesliga'
$newRows[4,3] = '2023-2024'
$newRows[4,4] = 45227.64583333334
$newRows[4,5] = 'Stuttgart'
$newRows[4,6] = 2
$newRows[4,7] = 'Hoffenheim'
$newRows[4,8] = 3
$newRows[4,9] = 1.77
$newRows[4,10] = '10/10/2023 14:02'
$newRows[4,11] = 1.67
$newRows[4,12] = '28/10/2023 15:27'
$newRows[4,13] = 4.11
$newRows[4,14] = '10/10/2023 14:02'
$newRows[4,15] = 4.43
$newRows[4,16] = '28/10/2023 15:29'
$newRows[4,17] = 4.38
$newRows[4,18] = '10/10/2023 14:02'
$newRows[4,19] = 4.84
$newRows[4,20] = '28/10/2023 15:29'
$newRows[4,21] = 'https://www.betexplorer.com/football/germany/bundesliga/vfb-stuttgart-hoffenheim/EuzDLv1F/'

$newRows[5,0] = 78
$newRows[5,1] = 'germany'
$newRows[5,2] = 'bundesliga'
$newRows[5,3] = '2023-2024'
$newRows[5,4] = 45227.64583333334
$newRows[5,5] = 'Werder Bremen'
$newRows[5,6] = 2
$newRows[5,7] = 'Union Berlin'
$newRows[5,8] = 0
$newRows[5,9] = 2.73
$newRows[5,10] = '10/10/2023 14:02'
$newRows[5,11] = 2.86
$newRows[5,12] = '28/10/2023 15:01'
$newRows[5,13] = 3.33
$newRows[5,14] = '10/10/2023 14:02'
$newRows[5,15] = 3.53
$newRows[5,16] = '28/10/2023 14:50'
$newRows[5,17] = 2.73
$newRows[5,18] = '10/10/2023 14:02'
$newRows[5,19] = 2.53
$newRows[5,20] = '28/10/2023 15:01'
$newRows[5,21] = 'https://www.betexplorer.com/football/germany/bundesliga/werder-bremen-union-berlin/4AhMbHg2/'

$newRows[6,0] = 79
$newRows[6,1] = 'germany'
$newRows[6,2] = 'bundesliga'
$newRows[6,3] = '2023-2024'
$newRows[6,4] = 45227.77083333334
$newRows[6,5] = 'RB Leipzig'
$newRows[6,6] = 6
$newRows[6,7] = 'FC Koln'
$newRows[6,8] = 0
$newRows[6,9] = 1.35
$newRows[6,10] = '10/10/2023 14:02'
$newRows[6,11] = 1.49
$newRows[6,12] = '28/10/2023 18:16'
$newRows[6,13] = 5.33
$newRows[6,14] = '10/10/2023 14:02'
$newRows[6,15] = 4.78
$newRows[6,16] = '28/10/2023 18:28'
$newRows[6,17] = 7.45
$newRows[6,18] = '10/10/2023 14:02'
$newRows[6,19] = 6.56
$newRows[6,20] = '28/10/2023 18:29'
$newRows[6,21] = 'https://www.betexplorer.com/football/germany/bundesliga/rb-leipzig-1-fc-koln/OhxtJouH/'

$newRows[7,0] = 80
$newRows[7,1] = 'germany'
$newRows[7,2] = 'bundesliga'
$newRows[7,3] = '2023-2024'
$newRows[7,4] = 45228.64583333334
$newRows[7,5] = 'Eintracht Frankfurt'
$newRows[7,6] = 3
$newRows[7,7] = 'Dortmund'
$newRows[7,8] = 3
$newRows[7,9] = 2.83
$newRows[7,10] = '10/10/2023 14:02'
$newRows[7,11] = 3.14
$newRows[7,12] = '29/10/2023 15:27'
$newRows[7,13] = 3.69
$newRows[7,14] = '10/10/2023 14:02'
$newRows[7,15] = 3.7
$newRows[7,16] = '29/10/2023 15:28'
$newRows[7,17] = 2.32
$newRows[7,18] = '10/10/2023 14:02'
$newRows[7,19] = 2.28
$newRows[7,20] = '29/10/2023 15:29'
$newRows[7,21] = 'https://www.betexplorer.com/football/germany/bundesliga/eintracht-frankfurt-dortmund/8raD0wPk/'

$newRows[8,0] = 81
$newRows[8,1] = 'germany'
$newRows[8,2] = 'bundesliga'
$newRows[8,3] = '2023-2024'
$newRows[8,4] = 45228.72916666666
$newRows[8,5] = 'Bayer Leverkusen'
$newRows[8,6] = 2
$newRows[8,7] = 'Freiburg'
$newRows[8,8] = 1
$newRows[8,9] = 1.45
$newRows[8,10] = '10/10/2023 14:02'
$newRows[8,11] = 1.24
$newRows[8,12] = '29/10/2023 17:28'
$newRows[8,13] = 4.87
$newRows[8,14] = '10/10/2023 14:02'
$newRows[8,15] = 6.67
$newRows[8,16] = '29/10/2023 17:29'
$newRows[8,17] = 7.04
$newRows[8,18] = '10/10/2023 14:02'
$newRows[8,19] = 12.56
$newRows[8,20] = '29/10/2023 17:29'
$newRows[8,21] = 'https://www.betexplorer.com/football/germany/bundesliga/bayer-leverkusen-freiburg/Ui091J9q/'

$ws.Range("A74:V82").Value2 = $newRows
